$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2: Santa Catarina
$ws.Range("A2").Value = "Santa Catarina"
Set-TextValue $ws.Range("C2") "01/10/2023"
$ws.Range("D2").Value = 96.8041974719771

# Row 3: Rondônia
$ws.Range("A3").Value = "Rondônia"
Set-TextValue $ws.Range("C3") "01/10/2023"
$ws.Range("D3").Value = 96.30512514898689

# Row 4: Mato Grosso
$ws.Range("A4").Value = "Mato Grosso"
Set-TextValue $ws.Range("C4") "01/10/2023"
$ws.Range("D4").Value = 96.13998970663921

# Row 5: Mato Grosso do Sul (name unchanged)
Set-TextValue $ws.Range("C5") "01/10/2023"
$ws.Range("D5").Value = 95.99733155436958

# Row 6: Paraná (name unchanged)
Set-TextValue $ws.Range("C6") "01/10/2023"
$ws.Range("D6").Value = 95.29524723955834

# Row 7: Rio Grande do Sul
$ws.Range("A7").Value = "Rio Grande do Sul"
Set-TextValue $ws.Range("C7") "01/10/2023"
$ws.Range("D7").Value = 94.8190658377172

# Row 8: Sergipe (name unchanged)
Set-TextValue $ws.Range("C8") "01/10/2023"
$ws.Range("D8").Value = 88.70214752567693
$ws.Range("E8").Value = "24º"

# Row 9: Nordeste (name unchanged)
Set-TextValue $ws.Range("C9") "01/10/2023"
$ws.Range("D9").Value = 89.56449309852451

# Row 10: Brasil (name unchanged)
Set-TextValue $ws.Range("C10") "01/10/2023"
$ws.Range("D10").Value = 92.59072488218143
